# Auto-generated market-price refresh for the Leve profitability workbook.
# For each affected (sheet, row) this sets the updated currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H,I,J,K,L,M,N) to the refreshed values pulled
# by the scheduled runner. Where a profit column is no longer applicable the cell
# is cleared (set to $null) instead of holding a stale number; one previously-blank
# profit cell (CUL!M120) gains a freshly computed value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2006.96
$ws.Range("J17").Value = 2006.96
$ws.Range("L17").Value = 6020.88
$ws.Range("N17").Value = -6356.88

$ws.Range("H33").Value = 173.07143
$ws.Range("I33").Value = 173.07143
$ws.Range("K33").Value = 173.07143
$ws.Range("M33").Value = 55.92857000000001

$ws.Range("H135").Value = 83334090
$ws.Range("I135").Value = 870.5
$ws.Range("J135").Value = 250000510
$ws.Range("K135").Value = 7834.5
$ws.Range("L135").Value = 2250004590
$ws.Range("M135").Value = -5299.5
$ws.Range("N135").Value = -2250009660


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 20451.334
$ws.Range("J24").Value = 20451.334
$ws.Range("L24").Value = 20451.334
$ws.Range("N24").Value = -21199.334

$ws.Range("H61").Value = 2777.111
$ws.Range("I61").Value = 2777.111
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2777.111
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2565.111
$ws.Range("N61").Value = $null

$ws.Range("H92").Value = 119999.8
$ws.Range("J92").Value = 119999.8
$ws.Range("L92").Value = 119999.8
$ws.Range("N92").Value = -124991.8

$ws.Range("H95").Value = 67217
$ws.Range("J95").Value = 67217
$ws.Range("L95").Value = 67217
$ws.Range("N95").Value = -72709

$ws.Range("H100").Value = 20451.334
$ws.Range("J100").Value = 20451.334
$ws.Range("L100").Value = 20451.334
$ws.Range("N100").Value = -22615.334

$ws.Range("H115").Value = 50684
$ws.Range("J115").Value = 50684
$ws.Range("L115").Value = 50684
$ws.Range("N115").Value = -53818

$ws.Range("H136").Value = 2777.111
$ws.Range("I136").Value = 2777.111
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8331.332999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5781.332999999999
$ws.Range("N136").Value = $null


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2849.6667
$ws.Range("J64").Value = 2831.875
$ws.Range("L64").Value = 2831.875
$ws.Range("N64").Value = -3281.875

$ws.Range("H67").Value = 2849.6667
$ws.Range("J67").Value = 2831.875
$ws.Range("L67").Value = 2831.875
$ws.Range("N67").Value = -4391.875

$ws.Range("H74").Value = 80879.5
$ws.Range("I74").Value = 80999
$ws.Range("J74").Value = 80760
$ws.Range("K74").Value = 80999
$ws.Range("L74").Value = 80760
$ws.Range("M74").Value = -80063
$ws.Range("N74").Value = -82632

$ws.Range("H77").Value = 80879.5
$ws.Range("I77").Value = 80999
$ws.Range("J77").Value = 80760
$ws.Range("K77").Value = 242997
$ws.Range("L77").Value = 242280
$ws.Range("M77").Value = -238317
$ws.Range("N77").Value = -251640

$ws.Range("H81").Value = 68390
$ws.Range("J81").Value = 68390
$ws.Range("L81").Value = 68390
$ws.Range("N81").Value = -70512

$ws.Range("H84").Value = 68390
$ws.Range("J84").Value = 68390
$ws.Range("L84").Value = 205170
$ws.Range("N84").Value = -215778


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 918.5
$ws.Range("I22").Value = 923.2
$ws.Range("J22").Value = 895
$ws.Range("K22").Value = 923.2
$ws.Range("L22").Value = 895
$ws.Range("M22").Value = -573.2
$ws.Range("N22").Value = -1595


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 228
$ws.Range("I47").Value = 80
$ws.Range("J47").Value = 302
$ws.Range("K47").Value = 240
$ws.Range("L47").Value = 906
$ws.Range("M47").Value = 191
$ws.Range("N47").Value = -1768

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null

$ws.Range("H68").Value = 3023.182
$ws.Range("I68").Value = 2059.2
$ws.Range("J68").Value = 3306.7058
$ws.Range("K68").Value = 6177.599999999999
$ws.Range("L68").Value = 9920.117400000001
$ws.Range("M68").Value = -5366.599999999999
$ws.Range("N68").Value = -11542.1174

$ws.Range("H71").Value = 3023.182
$ws.Range("I71").Value = 2059.2
$ws.Range("J71").Value = 3306.7058
$ws.Range("K71").Value = 18532.8
$ws.Range("L71").Value = 29760.3522
$ws.Range("M71").Value = -14476.8
$ws.Range("N71").Value = -37872.3522

$ws.Range("H82").Value = 500
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = $null

$ws.Range("H85").Value = 500
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = $null

$ws.Range("H115").Value = 3000
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").Value = $null

$ws.Range("H120").Value = 15000
$ws.Range("I120").Value = 15000
$ws.Range("K120").Value = 45000
$ws.Range("M120").Value = -40162


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 217666.17
$ws.Range("I3").Value = 201199.4
$ws.Range("J3").Value = 300000
$ws.Range("K3").Value = 201199.4
$ws.Range("L3").Value = 300000
$ws.Range("M3").Value = -201083.4
$ws.Range("N3").Value = -300232

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1138.6
$ws.Range("I22").Value = 1105.5714
$ws.Range("J22").Value = 1215.6666
$ws.Range("K22").Value = 1105.5714
$ws.Range("L22").Value = 1215.6666
$ws.Range("M22").Value = -810.5714
$ws.Range("N22").Value = -1805.6666

$ws.Range("H27").Value = 1138.6
$ws.Range("I27").Value = 1105.5714
$ws.Range("J27").Value = 1215.6666
$ws.Range("K27").Value = 1105.5714
$ws.Range("L27").Value = 1215.6666
$ws.Range("M27").Value = -998.5714
$ws.Range("N27").Value = -1429.6666

$ws.Range("H101").Value = 19709
$ws.Range("J101").Value = 19709
$ws.Range("L101").Value = 19709
$ws.Range("N101").Value = -26199

$ws.Range("H136").Value = 14275.125
$ws.Range("I136").Value = 26350.5
$ws.Range("J136").Value = 2199.75
$ws.Range("K136").Value = 79051.5
$ws.Range("L136").Value = 6599.25
$ws.Range("M136").Value = -76501.5
$ws.Range("N136").Value = -11699.25


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 80775.5
$ws.Range("J141").Value = 80775.5
$ws.Range("L141").Value = 80775.5
$ws.Range("N141").Value = -91135.5

